# Weekly update to the "Jengibre" (ginger) price sheet: a new weekly
# record is inserted at row 50 (pushing the existing rows 50-66 down to
# 51-67); all other columns for the new row repeat the constant values
# already used throughout the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new record by inserting a row above the current row 50.
$ws.Rows.Item(50).Insert()

# Populate the newly inserted row 50 with the new week's data.
$ws.Cells.Item(50, 1).Value  = 8
$ws.Cells.Item(50, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(50, 3).Value  = "Coquimbo"
$ws.Cells.Item(50, 4).Value  = 44837
$ws.Cells.Item(50, 5).Value  = 4
$ws.Cells.Item(50, 6).Value  = 100114007
$ws.Cells.Item(50, 7).Value  = "Jengibre"
$ws.Cells.Item(50, 8).Value  = "Sin especificar"
$ws.Cells.Item(50, 9).Value  = "Primera"
$ws.Cells.Item(50, 10).Value = 600
$ws.Cells.Item(50, 11).Value = 14000
$ws.Cells.Item(50, 12).Value = 15000
$ws.Cells.Item(50, 13).Value = 14500
$ws.Cells.Item(50, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(50, 15).Value = "Perú"
$ws.Cells.Item(50, 16).Value = 1115
$ws.Cells.Item(50, 17).Value = 13
$ws.Cells.Item(50, 18).Value = "Hortaliza"
